# Apply "break out stock.yaml completed" update to the "day" sheet:
#  1. Rows 697-726: the bsecode column (D) was entered as text in the
#     source feed; re-enter the same values as real numbers.
#  2. Append 20 new rows (727-746) pulled in from the 15/10/2024 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Fix column D (bsecode) typing for rows 697-726 -----------------
$bseFix = @(
    @{ Row=697; Bse=500290 },
    @{ Row=699; Bse=532500 },
    @{ Row=700; Bse=532977 },
    @{ Row=701; Bse=532538 },
    @{ Row=702; Bse=505200 },
    @{ Row=703; Bse=500495 },
    @{ Row=704; Bse=500114 },
    @{ Row=705; Bse=511218 },
    @{ Row=706; Bse=500520 },
    @{ Row=707; Bse=532343 },
    @{ Row=708; Bse=500325 },
    @{ Row=709; Bse=500790 },
    @{ Row=710; Bse=532478 },
    @{ Row=711; Bse=533309 },
    @{ Row=712; Bse=532978 },
    @{ Row=713; Bse=532921 },
    @{ Row=714; Bse=500271 },
    @{ Row=715; Bse=532868 },
    @{ Row=716; Bse=532733 },
    @{ Row=717; Bse=507685 },
    @{ Row=718; Bse=524208 },
    @{ Row=719; Bse=533278 },
    @{ Row=720; Bse=500875 },
    @{ Row=721; Bse=500400 },
    @{ Row=722; Bse=535755 },
    @{ Row=723; Bse=532720 },
    @{ Row=724; Bse=500103 },
    @{ Row=725; Bse=517334 },
    @{ Row=726; Bse=500113 }
)

foreach ($fix in $bseFix) {
    $ws.Cells.Item($fix.Row, 4).Value = $fix.Bse
}

# --- 2. Append the new rows 727-746 -------------------------------------
$newRows = @(
    @{ Row=727; Sr=1;  Code="APOLLOHOSP"; Name="Apollo Hospitals Enterprise Limited";            Bse="508869"; Chg=0.67;  Close=7140.15; Vol=322659;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=728; Sr=2;  Code="BAJFINANCE"; Name="Bajaj Finance Limited";                           Bse="500034"; Chg=-2.66; Close=7016.9;  Vol=1485567;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=729; Sr=3;  Code="LTIM";       Name="LTI Mindtree Ltd";                                Bse="540005"; Chg=0.19;  Close=6460.8;  Vol=597153;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=730; Sr=4;  Code="COLPAL";     Name="Colgate Palmolive (india) Limited";               Bse="500830"; Chg=-0.86; Close=3492.05; Vol=442238;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=731; Sr=5;  Code="PIDILITIND"; Name="Pidilite Industries Limited";                     Bse="500331"; Chg=1.03;  Close=3163.75; Vol=775805;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=732; Sr=6;  Code="SRF";        Name="Srf Limited";                                     Bse="503806"; Chg=0.35;  Close=2351.05; Vol=288559;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=733; Sr=7;  Code="ACC";        Name="Acc Limited";                                     Bse="500410"; Chg=-0.98; Close=2294.8;  Vol=356561;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=734; Sr=8;  Code="SUNPHARMA";  Name="Sun Pharmaceuticals Industries Limited";          Bse="524715"; Chg=-0.59; Close=1898.45; Vol=1383270;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=735; Sr=9;  Code="TECHM";      Name="Tech Mahindra Limited";                           Bse="532755"; Chg=-1.01; Close=1675.35; Vol=2503381;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=736; Sr=10; Code="BATAINDIA";  Name="Bata India Limited";                              Bse="500043"; Chg=2.47;  Close=1450.65; Vol=288039;   Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=737; Sr=11; Code="AMBUJACEM";  Name="Ambuja Cements Limited";                          Bse="500425"; Chg=0.25;  Close=590.35;  Vol=1805214;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=738; Sr=12; Code="VEDL";       Name="Vedanta Limited";                                 Bse="500295"; Chg=-1.86; Close=489.85;  Vol=7413480;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=739; Sr=13; Code="NTPC";       Name="Ntpc Limited";                                    Bse="532555"; Chg=0.49;  Close=426.6;   Vol=9365964;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=740; Sr=14; Code="HINDPETRO";  Name="Hindustan Petroleum Corporation Limited";         Bse="500104"; Chg=4.2;   Close=422.9;   Vol=16856414; Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=741; Sr=15; Code="BPCL";       Name="Bharat Petroleum Corporation Limited";            Bse="500547"; Chg=2.35;  Close=348.75;  Vol=17183945; Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=742; Sr=16; Code="GAIL";       Name="Gail (india) Limited";                            Bse="532155"; Chg=0.24;  Close=231.23;  Vol=16946401; Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=743; Sr=17; Code="IOC";        Name="Indian Oil Corporation Limited";                  Bse="530965"; Chg=1.49;  Close=167.93;  Vol=22551022; Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=744; Sr=18; Code="TATASTEEL";  Name="Tata Steel Limited";                               Bse="500470"; Chg=-1.7;  Close=155.63;  Vol=38964393; Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=745; Sr=19; Code="CANBK";      Name="Canara Bank";                                      Bse="532483"; Chg=-0.06; Close=104.43;  Vol=8407143;  Tf="day"; Dt="15/10/2024 11:36:35" },
    @{ Row=746; Sr=20; Code="IDFCFIRSTB"; Name="IDFC First Bank Ltd";                              Bse="539437"; Chg=-0.27; Close=72.74;   Vol=10389342; Tf="day"; Dt="15/10/2024 11:36:35" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Sr
    $ws.Cells.Item($r.Row, 2).Value = $r.Code
    $ws.Cells.Item($r.Row, 3).Value = $r.Name
    # bsecode stays text for the freshly appended rows (matches source feed) -
    # force text typing with a leading apostrophe so the numeric-looking
    # string isn't auto-converted to a number.
    $ws.Cells.Item($r.Row, 4).Value = "'" + $r.Bse
    $ws.Cells.Item($r.Row, 5).Value = $r.Chg
    $ws.Cells.Item($r.Row, 6).Value = $r.Close
    $ws.Cells.Item($r.Row, 7).Value = $r.Vol
    $ws.Cells.Item($r.Row, 8).Value = $r.Tf
    $ws.Cells.Item($r.Row, 9).Value = $r.Dt
}

Write-Host "Updated day sheet: fixed $($bseFix.Count) bsecode cells, appended $($newRows.Count) rows."
